$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new blank row at position 29 (pushes old row29+ down by one) ---
$ws.Rows("29:29").Insert()

# --- Column B width change ---
$ws.Columns("B:B").ColumnWidth = 16.7109375

# --- Row 10 ---
$ws.Range("B10").Value = "cate este cazul"

# --- Row 11 ---
$ws.Range("B11").Value = "cate este cazul"

# --- Row 13 ---
$ws.Range("B13").Value = "cate este cazul"

# --- Row 14 ---
$ws.Range("B14").Value = "cate este cazul"

# --- Row 15 ---
$ws.Range("C15").Value = "Sau Orice alta placa de control."

# --- Row 16 (new part, Rulment Axial) ---
$ws.Range("A16").Value = "Rulment Axial 51103"
$ws.Range("B16").NumberFormat = "d-mmm"
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("B16").Value = "1 sau 3"
$ws.Range("C16").Value = "3 în cazul în care se va implementa un pat incălzit (incă în lucru)"

# --- Row 17 ---
$ws.Range("A17").Value = "Sursa de Putere MW RSP 320"

# --- Row 18 ---
$ws.Range("A18").Value = "Bloc Incalzire"
$ws.Range("B18").Value = 1

# --- Row 19 ---
$ws.Range("A19").Value = "Senzori de Temperatura  NTC 100k / 10k(Pat*)"

# --- Row 20 ---
$ws.Range("A20").Value = "Senzori de Limita "
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = "Cate Unul pentru fiecare axa ( X și Z)"

# --- Row 21 ---
$ws.Range("A21").Value = "Motoare Pas Cu Pas Nema 17"
$ws.Range("B21").Value = 4

# --- Row 22 ---
$ws.Range("A22").Value = "Surub Filetat T8 - 300 mm"

# --- Row 23 ---
$ws.Range("A23").Value = "Cuplaj Surub Filetat"
$ws.Range("C23").ClearContents()

# --- Row 24 ---
$ws.Range("A24").Value = "Ecran de control"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = "Orice ecran compatibil cu Marlin Screen (adica aprope orice)."

# --- Row 25 ---
$ws.Range("A25").Value = "Cuplaj Curea de Transmisie"

# --- Row 26 ---
$ws.Range("A26").Value = "Rulment cu Flanja  6 x 12 x 4 mm"
$ws.Range("B26").Value = 2

# --- Row 27 ---
$ws.Range("A27").Value = "Piulita Trapeizoidala T8"
$ws.Range("B27").Value = 1

# --- Row 28 ---
$ws.Range("A28").Value = "Piese printate si cabluri"
$ws.Range("B28").Value = "Necesar"
$ws.Range("C28").ClearContents()

# --- Row 29 (new Miscellaneous row) ---
$ws.Range("A29").Value = "Miscellaneous"
$ws.Range("B29").Value = "-"
$ws.Range("C29").Value = "Unelte,conectori,zip-ties,cabluri USB, etc"

# --- Row 30 ---
$ws.Range("A30").Value = "Componente necesare pentru alte componente electronice"

# --- Row 31 ---
$ws.Range("A31").Value = "Driver Motor pas cu pas (TMC2208/9)"
$ws.Range("B31").Value = 4

# --- sheetView: topLeftCell + selection ---
$ws.Application.GoTo($ws.Range("C29"), $true)
$ws.Range("C29").Select()
